$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width updates (closest representable values; COM ColumnWidth
# quantizes to whole-pixel granularity, so 15.7109375 / 16.42578125
# round-trip to 15.666666666666666 / 16.5 once saved).
$ws.Columns.Item(1).ColumnWidth = 14.83
$ws.Columns.Item(2).ColumnWidth = 15.65

# Cell value updates
$ws.Range("A1").Value = -0.21540988069379097
$ws.Range("B1").Value = 0.21492789778510257
$ws.Range("A2").Value = -0.17310866161972083
$ws.Range("B2").Value = 0.17146933227789507
$ws.Range("A3").Value = -0.12175304177796242
$ws.Range("B3").Value = 0.12124745233834311
$ws.Range("A4").Value = -0.11324745246233547
$ws.Range("B4").Value = 0.11279250109313566
$ws.Range("A5").Value = -0.10979250116653638
$ws.Range("B5").Value = 0.10824233871341349
$ws.Range("A6").Value = -0.0097628029970806551
$ws.Range("B6").Value = 0.0097275166201562513
$ws.Range("A7").Value = 0.00027248320438699736
$ws.Range("B7").Value = -0.00027280733570611559
$ws.Range("A8").Value = 0.010272807160633501
$ws.Range("B8").Value = -0.01028510997936305
$ws.Range("A9").Value = 0.012285109901711166
$ws.Range("B9").Value = -0.012308800432854206
$ws.Range("A10").Value = 0.014308800358824314
$ws.Range("B10").Value = -0.014308882379561183
$ws.Range("A11").Value = -0.017862608351359555
$ws.Range("B11").Value = 0.017857487788940141
$ws.Range("A12").Value = -0.014357487881217601
$ws.Range("B12").Value = 0.014321311225079736
$ws.Range("A13").Value = -0.010821311320802884
$ws.Range("B13").Value = 0.010811571508273055
$ws.Range("A14").Value = -0.043516647216043758
$ws.Range("B14").Value = 0.043352544702243812
$ws.Range("A15").Value = -0.042352544774529655
$ws.Range("B15").Value = 0.042207247127977965
$ws.Range("A16").Value = -0.0060342281310097334
$ws.Range("B16").Value = 0.006003052153217503
$ws.Range("A17").Value = -0.004003052245329819
$ws.Range("B17").Value = 0.0039999998838782247
$ws.Range("A18").Value = -0.016100204915328931
$ws.Range("B18").Value = 0.01609070496360232
$ws.Range("A19").Value = -0.012090705014023762
$ws.Range("B19").Value = 0.012015717679518811
$ws.Range("A20").Value = -0.0080157177337145669
$ws.Range("B20").Value = 0.0080055663443410197
$ws.Range("A21").Value = -0.0040055663991171997
$ws.Range("B21").Value = 0.0039999999447157819
$ws.Range("A22").Value = -0.045716290568046603
$ws.Range("B22").Value = 0.045502233234989475
$ws.Range("A23").Value = -0.040502233319133829
$ws.Range("B23").Value = 0.040099453929434148
$ws.Range("A24").Value = -0.020099454196877531
$ws.Range("B24").Value = 0.019999999728751661
$ws.Range("A25").Value = -0.0972784067426673
$ws.Range("B25").Value = 0.097154209896215704
$ws.Range("A26").Value = -0.09465420998753693
$ws.Range("B26").Value = 0.094492768232800728
$ws.Range("A27").Value = -0.091992768329917762
$ws.Range("B27").Value = 0.09102656806277265
$ws.Range("A28").Value = -0.089026568177160037
$ws.Range("B28").Value = 0.088365328556741396
$ws.Range("A29").Value = -0.081365328746876742
$ws.Range("B29").Value = 0.081174340354320584
$ws.Range("A30").Value = -0.021174341178350264
$ws.Range("B30").Value = 0.021022018130464115
$ws.Range("A31").Value = -0.014022018334626907
$ws.Range("B31").Value = 0.014000627306115732
$ws.Range("A32").Value = -0.0040006275461781371
$ws.Range("B32").Value = 0.0039999998305439988
